$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 0) Dashboard_Tests: leftover selection change from editing (C15 -> C12).
# ------------------------------------------------------------------
$dashboardWs = $wb.Worksheets.Item("Dashboard_Tests")
$dashboardWs.Range("C12").Select()

# ------------------------------------------------------------------
# 1) Sales_Tests: expand TC_SALE_02 ("Verify List Consistency") into
#    "Verify List Consistency & Invoice" with extra steps for viewing
#    the invoice, verifying the URL/heading, and navigating back.
# ------------------------------------------------------------------
$salesWs = $wb.Worksheets.Item("Sales_Tests")

# Insert 4 new blank rows right after row 4 (the TC_SALE_02 header row)
# to make room for the new invoice-verification steps.
$salesWs.Rows.Item(5).Resize(4).Insert()

# Update the test-case description on row 4.
$salesWs.Range("B4").Value = "Verify List Consistency & Invoice"

# Fill in the newly inserted steps (rows 5-8). The shared-string table
# records brand-new text in the order it is first assigned, so these
# are written out-of-visual-order to line up with the authored file.
$salesWs.Range("C6").Value = '3.Verify URL contains "invoice"'
$salesWs.Range("C7").Value = '4.Verify text "NESTO SUPERMARKET" at "//h2"'
$salesWs.Range("C5").Value = '2.Click on "Show Invoice" at "(//a[contains(@href, ''/invoice/'')])[1]"'
$salesWs.Range("C8").Value = '5.Click "Back" at "//a[text()=''Back'']"'

# This sheet is no longer the active tab.
$salesWs.Range("C1").Select()

# ------------------------------------------------------------------
# 2) Add a new "Logout_Tests" sheet at the end of the workbook with a
#    single logout test case.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$logoutWs = $wb.Worksheets.Add($null, $lastSheet)
$logoutWs.Name = "Logout_Tests"

$logoutWs.Columns.Item(1).ColumnWidth = 17.83
$logoutWs.Columns.Item(2).ColumnWidth = 25.72
$logoutWs.Columns.Item(3).ColumnWidth = 59.28

$logoutWs.Range("A1").Value = "Test Case ID(s)"
$logoutWs.Range("B1").Value = "Test Case Description"
$logoutWs.Range("C1").Value = "Test Steps"
$logoutWs.Range("A1:C1").Interior.Color = 5287936

$logoutWs.Range("A2").Value = "TC_LOGOUT_01"
$logoutWs.Range("B2").Value = "Verify Logout Functionality"
$logoutWs.Range("C2").Value = '1.Click on "Dashboard Link" at "//a[contains(@href, ''/dashboard'')]"'
$logoutWs.Range("C3").Value = '2.Click on "Logout Button" at "//button[contains(@class, ''btn-logout'')]"'
$logoutWs.Range("C4").Value = '3.Verify URL contains "login"'

$logoutWs.Range("C1").Select()

# Logout_Tests is now the active sheet/tab.
$logoutWs.Activate()
